$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = [double]"2"
$ws.Cells.Item(2, 2).Value = [double]"0.163287555206"
$ws.Cells.Item(2, 3).Value = [double]"0.38149406496338573"
$ws.Cells.Item(2, 4).Value = [double]"4.2039687351188987E-2"
$ws.Cells.Item(2, 5).Value = [double]"2.9933333333333336"
$ws.Cells.Item(2, 6).Value = [double]"1.1547005383792526E-2"
$ws.Cells.Item(3, 1).Value = [double]"2.1"
$ws.Cells.Item(3, 2).Value = [double]"0.207831077394"
$ws.Cells.Item(3, 3).Value = [double]"0.54079969318947774"
$ws.Cells.Item(3, 4).Value = [double]"1.1307961032332713E-2"
$ws.Cells.Item(3, 5).Value = [double]"3.0033333333333334"
$ws.Cells.Item(3, 6).Value = [double]"5.7735026918961348E-3"
$ws.Cells.Item(4, 1).Value = [double]"2.2000000000000002"
$ws.Cells.Item(4, 2).Value = [double]"0.254589961024"
$ws.Cells.Item(4, 3).Value = [double]"0.54114371551005247"
$ws.Cells.Item(4, 4).Value = [double]"4.4944264440218446E-3"
$ws.Cells.Item(4, 5).Value = [double]"3.0033333333333334"
$ws.Cells.Item(4, 6).Value = [double]"8.3333333333332742E-3"
$ws.Cells.Item(5, 1).Value = [double]"2.2999999999999998"
$ws.Cells.Item(5, 2).Value = [double]"0.30112915976900001"
$ws.Cells.Item(5, 3).Value = [double]"0.54503464293473003"
$ws.Cells.Item(5, 4).Value = [double]"2.062412086264484E-2"
$ws.Cells.Item(5, 5).Value = [double]"2.9933333333333336"
$ws.Cells.Item(5, 6).Value = [double]"1.1547005383792526E-2"
$ws.Cells.Item(6, 1).Value = [double]"2.4"
$ws.Cells.Item(6, 2).Value = [double]"0.345328911889"
$ws.Cells.Item(6, 3).Value = [double]"0.53893602092389647"
$ws.Cells.Item(6, 4).Value = [double]"1.7379476258423448E-3"
$ws.Cells.Item(6, 5).Value = [double]"2.9933333333333336"
$ws.Cells.Item(6, 6).Value = [double]"1.1547005383792526E-2"
$ws.Cells.Item(7, 1).Value = [double]"2.5"
$ws.Cells.Item(7, 2).Value = [double]"0.38567879271299998"
$ws.Cells.Item(7, 3).Value = [double]"0.43834237035580376"
$ws.Cells.Item(7, 4).Value = [double]"1.7670864210681043E-2"
$ws.Cells.Item(7, 5).Value = [double]"5.0233333333333334"
$ws.Cells.Item(7, 6).Value = [double]"4.6188021535170098E-2"
$ws.Cells.Item(8, 1).Value = [double]"2.6"
$ws.Cells.Item(8, 2).Value = [double]"0.42134668099900002"
$ws.Cells.Item(8, 3).Value = [double]"0.54071101632868901"
$ws.Cells.Item(8, 4).Value = [double]"4.7569624366485246E-3"
$ws.Cells.Item(8, 5).Value = [double]"4.99"
$ws.Cells.Item(8, 6).Value = [double]"2.6457513110645845E-2"
$ws.Cells.Item(9, 1).Value = [double]"2.7"
$ws.Cells.Item(9, 2).Value = [double]"0.45208279949800001"
$ws.Cells.Item(9, 3).Value = [double]"0.5542746649141852"
$ws.Cells.Item(9, 4).Value = [double]"6.2283865550030989E-3"
$ws.Cells.Item(9, 5).Value = [double]"4.9933333333333332"
$ws.Cells.Item(9, 6).Value = [double]"2.0816659994661382E-2"
$ws.Cells.Item(10, 1).Value = [double]"2.8"
$ws.Cells.Item(10, 2).Value = [double]"0.47805208623700002"
$ws.Cells.Item(10, 3).Value = [double]"0.55710007035718179"
$ws.Cells.Item(10, 4).Value = [double]"6.8005762380074983E-3"
$ws.Cells.Item(10, 5).Value = [double]"5.0366666666666662"
$ws.Cells.Item(10, 6).Value = [double]"1.5275252316519626E-2"
$ws.Cells.Item(11, 1).Value = [double]"2.9000000000099999"
$ws.Cells.Item(11, 2).Value = [double]"0.49966675215799999"
$ws.Cells.Item(11, 3).Value = [double]"0.5764148930595675"
$ws.Cells.Item(11, 4).Value = [double]"3.1099084834229291E-3"
$ws.Cells.Item(11, 5).Value = [double]"5.0133333333333328"
$ws.Cells.Item(11, 6).Value = [double]"4.1633319989322265E-2"
$ws.Cells.Item(12, 1).Value = [double]"3.00000000001"
$ws.Cells.Item(12, 2).Value = [double]"0.51745431077500004"
$ws.Cells.Item(12, 3).Value = [double]"0.469249011180983"
$ws.Cells.Item(12, 4).Value = [double]"2.2708007255614518E-2"
$ws.Cells.Item(12, 5).Value = [double]"6.9733333333333336"
$ws.Cells.Item(12, 6).Value = [double]"3.2145502536643E-2"
$ws.Cells.Item(13, 1).Value = [double]"3.1"
$ws.Cells.Item(13, 2).Value = [double]"0.53196825941799997"
$ws.Cells.Item(13, 3).Value = [double]"0.57150875879833707"
$ws.Cells.Item(13, 4).Value = [double]"3.8679315018278803E-3"
$ws.Cells.Item(13, 5).Value = [double]"6.9833333333333334"
$ws.Cells.Item(13, 6).Value = [double]"3.0550504633038766E-2"
$ws.Cells.Item(14, 1).Value = [double]"3.2"
$ws.Cells.Item(14, 2).Value = [double]"0.54373576958799996"
$ws.Cells.Item(14, 3).Value = [double]"0.58323433537480052"
$ws.Cells.Item(14, 4).Value = [double]"5.3911615325178009E-3"
$ws.Cells.Item(14, 5).Value = [double]"6.9733333333333336"
$ws.Cells.Item(14, 6).Value = [double]"2.5385910352879595E-2"
$ws.Cells.Item(15, 1).Value = [double]"3.3"
$ws.Cells.Item(15, 2).Value = [double]"0.55323088734899994"
$ws.Cells.Item(15, 3).Value = [double]"0.59487480074657484"
$ws.Cells.Item(15, 4).Value = [double]"2.1472861288731531E-2"
$ws.Cells.Item(15, 5).Value = [double]"7.0166666666666666"
$ws.Cells.Item(15, 6).Value = [double]"2.5166114784235707E-2"
$ws.Cells.Item(16, 1).Value = [double]"3.4"
$ws.Cells.Item(16, 2).Value = [double]"0.56086476759199999"
$ws.Cells.Item(16, 3).Value = [double]"0.59614898791591653"
$ws.Cells.Item(16, 4).Value = [double]"3.1328926411537312E-3"
$ws.Cells.Item(16, 5).Value = [double]"6.9833333333333334"
$ws.Cells.Item(16, 6).Value = [double]"3.2145502536643007E-2"
$ws.Cells.Item(17, 1).Value = [double]"3.5"
$ws.Cells.Item(17, 2).Value = [double]"0.56698544995699995"
$ws.Cells.Item(17, 3).Value = [double]"0.5178425389738156"
$ws.Cells.Item(17, 4).Value = [double]"3.6599482021059691E-2"
$ws.Cells.Item(17, 5).Value = [double]"8.99"
$ws.Cells.Item(17, 6).Value = [double]"1.7320508075688402E-2"
$ws.Cells.Item(18, 1).Value = [double]"3.6"
$ws.Cells.Item(18, 2).Value = [double]"0.57188264455600002"
$ws.Cells.Item(18, 3).Value = [double]"0.62208412840144678"
$ws.Cells.Item(18, 4).Value = [double]"1.3625560510124042E-2"
$ws.Cells.Item(18, 5).Value = [double]"8.9866666666666664"
$ws.Cells.Item(18, 6).Value = [double]"1.154700538379227E-2"
$ws.Cells.Item(19, 1).Value = [double]"3.7"
$ws.Cells.Item(19, 2).Value = [double]"0.57579465540700003"
$ws.Cells.Item(19, 3).Value = [double]"0.63300654045813587"
$ws.Cells.Item(19, 4).Value = [double]"1.5685037377726563E-3"
$ws.Cells.Item(19, 5).Value = [double]"8.9800000000000022"
$ws.Cells.Item(19, 6).Value = [double]"9.9999999999997868E-3"
$ws.Cells.Item(20, 1).Value = [double]"3.8"
$ws.Cells.Item(20, 2).Value = [double]"0.57891582620699999"
$ws.Cells.Item(20, 3).Value = [double]"0.60939630849511151"
$ws.Cells.Item(20, 4).Value = [double]"1.8529052300006489E-2"
$ws.Cells.Item(20, 5).Value = [double]"8.9966666666666661"
$ws.Cells.Item(20, 6).Value = [double]"2.5166114784235295E-2"
$ws.Cells.Item(21, 1).Value = [double]"3.9000000000999999"
$ws.Cells.Item(21, 2).Value = [double]"0.581403656913"
$ws.Cells.Item(21, 3).Value = [double]"0.60923084026543728"
$ws.Cells.Item(21, 4).Value = [double]"1.5296611248176354E-2"
$ws.Cells.Item(21, 5).Value = [double]"8.99"
$ws.Cells.Item(21, 6).Value = [double]"9.9999999999997868E-3"
$ws.Cells.Item(22, 1).Value = [double]"4.00000000001"
$ws.Cells.Item(22, 2).Value = [double]"0.58338519306000003"
$ws.Cells.Item(23, 1).Value = [double]"4.0999999999999996"
$ws.Cells.Item(23, 2).Value = [double]"0.58496255613399994"
$ws.Cells.Item(24, 1).Value = [double]"4.2"
$ws.Cells.Item(24, 2).Value = [double]"0.58621761515299997"
$ws.Cells.Item(25, 1).Value = [double]"4.3"
$ws.Cells.Item(25, 2).Value = [double]"0.58721586981100005"
$ws.Cells.Item(26, 1).Value = [double]"4.4000000000000004"
$ws.Cells.Item(26, 2).Value = [double]"0.58800964140800005"
$ws.Cells.Item(27, 1).Value = [double]"4.5"
$ws.Cells.Item(27, 2).Value = [double]"0.58864067366999995"
$ws.Cells.Item(28, 1).Value = [double]"4.5999999999999996"
$ws.Cells.Item(28, 2).Value = [double]"0.58914223957800005"
$ws.Cells.Item(29, 1).Value = [double]"4.7"
$ws.Cells.Item(29, 2).Value = [double]"0.58954084055"
$ws.Cells.Item(30, 1).Value = [double]"4.8"
$ws.Cells.Item(30, 2).Value = [double]"0.58985757212400003"
$ws.Cells.Item(31, 1).Value = [double]"4.9000000000000004"
$ws.Cells.Item(31, 2).Value = [double]"0.59010921867199995"
$ws.Cells.Item(32, 1).Value = [double]"5"
$ws.Cells.Item(32, 2).Value = [double]"0.59030912926000001"
$ws.Cells.Item(33, 1).Value = [double]"5.0999999999999996"
$ws.Cells.Item(33, 2).Value = [double]"0.59046791720500003"
$ws.Cells.Item(34, 1).Value = [double]"5.2"
$ws.Cells.Item(34, 2).Value = [double]"0.59059401819799995"
$ws.Cells.Item(35, 1).Value = [double]"5.3"
$ws.Cells.Item(35, 2).Value = [double]"0.590694135111"
$ws.Cells.Item(36, 1).Value = [double]"5.4"
$ws.Cells.Item(36, 2).Value = [double]"0.59077359224000003"
$ws.Cells.Item(37, 1).Value = [double]"5.5"
$ws.Cells.Item(37, 2).Value = [double]"0.59083661724600001"
$ws.Cells.Item(38, 1).Value = [double]"5.6"
$ws.Cells.Item(38, 2).Value = [double]"0.59088656543899998"
$ws.Cells.Item(39, 1).Value = [double]"5.7"
$ws.Cells.Item(39, 2).Value = [double]"0.59092609815899999"
$ws.Cells.Item(40, 1).Value = [double]"5.8"
$ws.Cells.Item(40, 2).Value = [double]"0.59095732465600004"
$ws.Cells.Item(41, 1).Value = [double]"5.9"
$ws.Cells.Item(41, 2).Value = [double]"0.59098191497300001"
$ws.Cells.Item(42, 1).Value = [double]"5.99999999988"
$ws.Cells.Item(42, 2).Value = [double]"0.59100121608699996"
$ws.Cells.Item(43, 1).Value = [double]"6.1"
$ws.Cells.Item(43, 2).Value = [double]"0.59101619270700001"
$ws.Cells.Item(44, 1).Value = [double]"6.1999999999500002"
$ws.Cells.Item(44, 2).Value = [double]"0.59102774675900005"
$ws.Cells.Item(45, 1).Value = [double]"6.3000000000499998"
$ws.Cells.Item(45, 2).Value = [double]"0.59103652925799999"
$ws.Cells.Item(46, 1).Value = [double]"6.4"
$ws.Cells.Item(46, 2).Value = [double]"0.59104297026599995"
$ws.Cells.Item(47, 1).Value = [double]"6.5"
$ws.Cells.Item(47, 2).Value = [double]"0.59104755907600004"
$ws.Cells.Item(48, 1).Value = [double]"6.6"
$ws.Cells.Item(48, 2).Value = [double]"0.59105059129799997"
$ws.Cells.Item(49, 1).Value = [double]"6.7"
$ws.Cells.Item(49, 2).Value = [double]"0.59105232506000005"
$ws.Cells.Item(50, 1).Value = [double]"6.8000000000499998"
$ws.Cells.Item(50, 2).Value = [double]"0.59105296598000001"
$ws.Cells.Item(51, 1).Value = [double]"6.9"
$ws.Cells.Item(51, 2).Value = [double]"0.59105267481299995"
$ws.Cells.Item(52, 1).Value = [double]"7"
$ws.Cells.Item(52, 2).Value = [double]"0.59105157066400005"
$ws.Cells.Item(53, 1).Value = [double]"7.1"
$ws.Cells.Item(53, 2).Value = [double]"0.59104973066400002"
$ws.Cells.Item(54, 1).Value = [double]"7.2"
$ws.Cells.Item(54, 2).Value = [double]"0.59104718721899996"
$ws.Cells.Item(55, 1).Value = [double]"7.3"
$ws.Cells.Item(55, 2).Value = [double]"0.59104392389799998"
$ws.Cells.Item(56, 1).Value = [double]"7.4"
$ws.Cells.Item(56, 2).Value = [double]"0.59103987064100005"
$ws.Cells.Item(57, 1).Value = [double]"7.5"
$ws.Cells.Item(57, 2).Value = [double]"0.59103489841800005"
$ws.Cells.Item(58, 1).Value = [double]"7.6"
$ws.Cells.Item(58, 2).Value = [double]"0.59102881288499998"
$ws.Cells.Item(59, 1).Value = [double]"7.7"
$ws.Cells.Item(59, 2).Value = [double]"0.59102128811599997"
$ws.Cells.Item(60, 1).Value = [double]"7.8"
$ws.Cells.Item(60, 2).Value = [double]"0.59101207135099998"
$ws.Cells.Item(61, 1).Value = [double]"7.9"
$ws.Cells.Item(61, 2).Value = [double]"0.59100066580800004"
$ws.Cells.Item(62, 1).Value = [double]"8"
$ws.Cells.Item(62, 2).Value = [double]"0.59098650146800003"
$ws.Cells.Item(63, 1).Value = [double]"8.1"
$ws.Cells.Item(63, 2).Value = [double]"0.59096886561799999"
$ws.Cells.Item(64, 1).Value = [double]"8.1999999999999993"
$ws.Cells.Item(64, 2).Value = [double]"0.590946872746"
$ws.Cells.Item(65, 1).Value = [double]"8.3000000000000007"
$ws.Cells.Item(65, 2).Value = [double]"0.59091942751299997"
$ws.Cells.Item(66, 1).Value = [double]"8.4"
$ws.Cells.Item(66, 2).Value = [double]"0.59088518209899998"
$ws.Cells.Item(67, 1).Value = [double]"8.5"
$ws.Cells.Item(67, 2).Value = [double]"0.59084248684200003"
$ws.Cells.Item(68, 1).Value = [double]"8.6"
$ws.Cells.Item(68, 2).Value = [double]"0.59078933703900005"
$ws.Cells.Item(69, 1).Value = [double]"8.6999999999999993"
$ws.Cells.Item(69, 2).Value = [double]"0.590723317784"
$ws.Cells.Item(70, 1).Value = [double]"8.8000000000000007"
$ws.Cells.Item(70, 2).Value = [double]"0.59064155441099997"
$ws.Cells.Item(71, 1).Value = [double]"8.9"
$ws.Cells.Item(71, 2).Value = [double]"0.59054067504300001"
$ws.Cells.Item(72, 1).Value = [double]"9"
$ws.Cells.Item(72, 2).Value = [double]"0.590416800962"
$ws.Cells.Item(73, 1).Value = [double]"9.1"
$ws.Cells.Item(73, 2).Value = [double]"0.59026558586599998"
$ws.Cells.Item(74, 1).Value = [double]"9.1999999997999993"
$ws.Cells.Item(74, 2).Value = [double]"0.59008232452800002"
$ws.Cells.Item(75, 1).Value = [double]"9.3000000000000007"
$ws.Cells.Item(75, 2).Value = [double]"0.58986216506300004"
$ws.Cells.Item(76, 1).Value = [double]"9.4"
$ws.Cells.Item(76, 2).Value = [double]"0.58960044901800002"
$ws.Cells.Item(77, 1).Value = [double]"9.5"
$ws.Cells.Item(77, 2).Value = [double]"0.58929319596999996"
$ws.Cells.Item(78, 1).Value = [double]"9.6"
$ws.Cells.Item(78, 2).Value = [double]"0.58893772149300005"
$ws.Cells.Item(79, 1).Value = [double]"9.6999999999999993"
$ws.Cells.Item(79, 2).Value = [double]"0.58853333355899995"
$ws.Cells.Item(80, 1).Value = [double]"9.8000000000499998"
$ws.Cells.Item(80, 2).Value = [double]"0.58808199674399997"
$ws.Cells.Item(81, 1).Value = [double]"9.9"
$ws.Cells.Item(81, 2).Value = [double]"0.58758880285500004"
$ws.Cells.Item(82, 1).Value = [double]"9.9999999999"
$ws.Cells.Item(82, 2).Value = [double]"0.58706206852800003"

$ws.Range("M14").Select() | Out-Null
